# Update the "your article inside holding" label (RU -> UA) and move the
# active selection to B1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Артикул всередині холдингу"

$ws.Range("B1").Select()
